$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "29.049.63"
Set-TextValue "E2" "  -0.73%  "
Set-TextValue "D3" "1.828.91"
Set-TextValue "E3" "  -0.71%  "
Set-TextValue "D4" "0.9993"
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "241.56"
Set-TextValue "E5" "  +0.27%  "
Set-TextValue "D6" "0.6302"
Set-TextValue "E6" "  -5.94%  "
Set-TextValue "D8" "44.66"
Set-TextValue "E8" "  +6.05%  "
Set-TextValue "D9" "0.2932"
Set-TextValue "E9" "  -0.12%  "
Set-TextValue "D10" "0.07331"
Set-TextValue "E10" "  -1.11%  "
Set-TextValue "D11" "22.87"
Set-TextValue "E11" "  +0.14%  "
Set-TextValue "D12" "0.07675"
Set-TextValue "E12" "  -0.47%  "
Set-TextValue "D13" "1.827.06"
Set-TextValue "E13" "  +0.06%  "
Set-TextValue "E14" "  -0.35%  "
Set-TextValue "D15" "0.6626"
Set-TextValue "E15" "  -1.18%  "
Set-TextValue "D16" "82.02"
Set-TextValue "E16" "  -4.50%  "
Set-TextValue "D17" "6.055"
Set-TextValue "E17" "  -1.45%  "
Set-TextValue "D18" "0.000008657"
Set-TextValue "D19" "29.037.77"
Set-TextValue "E19" "  -0.72%  "
Set-TextValue "D20" "2.080.34"
Set-TextValue "E20" "  +0.01%  "
Set-TextValue "D21" "12.39"
Set-TextValue "E21" "  -0.98%  "
Set-TextValue "D22" "223.83"
Set-TextValue "E22" "  -1.75%  "
Set-TextValue "D23" "1.000"
Set-TextValue "E23" "  -0.02%  "
Set-TextValue "D24" "7.126"
Set-TextValue "E24" "  -0.32%  "
Set-TextValue "E25" "  +0.00%  "
Set-TextValue "D26" "157.94"
Set-TextValue "E26" "  -1.87%  "
Set-TextValue "D27" "8.446"
Set-TextValue "E27" "  -2.88%  "
Set-TextValue "D28" "0.1366"
Set-TextValue "E28" "  -2.50%  "
Set-TextValue "E29" "  -0.85%  "
Set-TextValue "D30" "1.507"
Set-TextValue "E30" "  -0.37%  "
Set-TextValue "D31" "4.086"
Set-TextValue "E31" "  -1.53%  "
Set-TextValue "B32" "InternetComputer(DFINITY)"
Set-TextValue "C32" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D32" "4.018"
Set-TextValue "E32" "  -1.17%  "
Set-TextValue "B33" "Toncoin"
Set-TextValue "C33" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D33" "1.201"
Set-TextValue "E33" "  +0.46%  "
Set-TextValue "D34" "0.05299"
Set-TextValue "E34" "  -0.09%  "
Set-TextValue "D35" "0.7388"
Set-TextValue "E35" "  -1.68%  "
Set-TextValue "D36" "1.827"
Set-TextValue "E36" "  -2.62%  "
Set-TextValue "E37" "  +1.31%  "
Set-TextValue "E38" "  -1.08%  "
Set-TextValue "D39" "1.294.99"
Set-TextValue "E39" "  -1.80%  "
Set-TextValue "B40" "VeChain"
Set-TextValue "C40" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D40" "0.01781"
Set-TextValue "E40" "  -1.30%  "
Set-TextValue "B41" "MXToken"
Set-TextValue "C41" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D41" "2.739"
Set-TextValue "E41" "  +0.50%  "
Set-TextValue "D42" "6.318"
Set-TextValue "E42" "  +5.90%  "
Set-TextValue "D43" "0.8947"
Set-TextValue "E43" "  -2.74%  "
Set-TextValue "D44" "0.9997"
Set-TextValue "E44" "  -0.26%  "
Set-TextValue "D45" "102.59"
Set-TextValue "E45" "  +0.63%  "
Set-TextValue "B46" "RocketPoolETH"
Set-TextValue "C46" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D46" "1.978.20"
Set-TextValue "E46" "  -0.11%  "
Set-TextValue "B47" "BabyDogeCoin"
Set-TextValue "C47" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D47" "0.00000000124"
Set-TextValue "E47" "  +2.55%  "
Set-TextValue "D48" "0.5139"
Set-TextValue "E48" "  -0.50%  "
Set-TextValue "D49" "64.25"
Set-TextValue "E49" "  +0.71%  "
Set-TextValue "D50" "1.729"
Set-TextValue "E50" "  -2.51%  "
